$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark from the title paragraph.
#    (It gets re-added later, at the end of the "Project:" paragraph.)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) "Date:  " (two trailing spaces) -> "Date: " (one trailing space)
$d.Content.Find.Execute("Date:  ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Date: ", 2) | Out-Null

# 3) Underlined "DATE" run absorbs the following plain space run,
#    becoming a single underlined run reading "`DATE~".
$d.Content.Find.Execute("DATE ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "``DATE~", 2) | Out-Null

# 4) "Project: SHORT " -> four runs: "Project: ", "`", "SHORT", "~"
#    Locate the paragraph fresh (offsets may have shifted above).
$projStart = -1
$projEnd = -1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^Project:") {
        $projStart = $p.Range.Start
        $projEnd = $p.Range.End
    }
}

# "Project: " is 9 chars, so SHORT starts at $projStart + 9.
$shortStart = $projStart + 9
$shortEnd = $shortStart + 5          # "SHORT" is 5 chars
$spaceStart = $shortEnd              # trailing space right after SHORT

# Replace the trailing space with a tilde first (keeps offsets stable).
$rSpace = $d.Range($spaceStart, $spaceStart + 1)
$rSpace.Text = "~"

# Insert a backtick immediately before "SHORT".
$rIns = $d.Range($shortStart, $shortStart)
$rIns.InsertAfter("``")

# Force the backtick and tilde to stay in their own runs (distinct from
# their neighbours) by touching - then resetting - a formatting toggle.
# This must be the *last* thing done to each of those two characters so
# the run boundaries are not re-merged by a later paragraph edit.
$rBacktick = $d.Range($shortStart, $shortStart + 1)
$rBacktick.Bold = 1
$rBacktick.Bold = 0

$rTilde = $d.Range($spaceStart + 1, $spaceStart + 2)
$rTilde.Bold = 1
$rTilde.Bold = 0

# 5) Re-add the "_GoBack" bookmark at the very end of the "Project:"
#    paragraph (right before its paragraph mark). A temporary character is
#    inserted so the bookmark can anchor to a real (non-collapsed) range;
#    it is then removed, leaving the bookmark collapsed at that position.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^Project:") {
        $insertPos = $p.Range.End - 1
        $r = $d.Range($insertPos, $insertPos)
        $r.InsertAfter("X")
        $d.Bookmarks.Add("_GoBack", $r) | Out-Null
        $bm = $d.Bookmarks("_GoBack")
        $bm.Range.Text = ""
        break
    }
}
